# "Actualizar" refresh: the availability-check timestamp in column D is
# re-stamped with the current run time, and the previous run's timestamps
# shift down into the next older time-slot block.
#
#   rows 2-15  (slot 1, newest)  -> new "now" timestamp
#   rows 16-29 (slot 2)          -> gets the value that slot 1 used to have
#   rows 30-43 (slot 3, oldest)  -> gets the value that slot 2 used to have

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44233.554264075
$ws.Range("D16:D29").Value = 44233.53308902778
$ws.Range("D30:D43").Value = 44233.51190966435
